$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'task_id'
$ws.Range('B1').Value = 'prompt'
$ws.Range('C1').Value = 'canonical_solution'
$ws.Range('D1').Value = 'result'

$ws.Range('A2').Value = 'test_sql/0'
$ws.Range('B2').Value = 'How many products do we have? | products : product_id, country, price'
$ws.Range('C2').Value = '[''SELECT COUNT(*) FROM PRODUCTS'']'
$ws.Range('D2').Value = '[False]'

$ws.Range('A3').Value = 'test_sql/1'
$ws.Range('B3').Value = 'What is the total number of orders? | orders : order_id, description'
$ws.Range('C3').Value = '[''SELECT COUNT(*) FROM ORDERS'']'
$ws.Range('D3').Value = '[False]'

$ws.Range('A4').Value = 'test_sql/2'
$ws.Range('B4').Value = 'What is the average , minimum , and maximum price of all Spanish products? | products : product_id, country, price'
$ws.Range('C4').Value = '["SELECT AVG(PRICE), MIN(PRICE), MAX(PRICE) FROM PRODUCTS WHERE COUNTRY = ''SPAIN''"]'
$ws.Range('D4').Value = '[False]'

$ws.Range('A5').Value = 'test_sql/3'
$ws.Range('B5').Value = 'Show all countries and the number of products in each country | products'
$ws.Range('C5').Value = '[''SELECT COUNTRY, COUNT(*) FROM PRODUCTS GROUP BY COUNTRY'']'
$ws.Range('D5').Value = '[False]'

$ws.Range('A6').Value = 'test_sql/4'
$ws.Range('B6').Value = 'How many sales are there in store STORE1? | sales: sale_id, product_id, customer_id, quantity, store'
$ws.Range('C6').Value = '["SELECT STORE, COUNT(*) FROM SALES WHERE BY STORE = ''STORE1''"]'
$ws.Range('D6').Value = '[False]'

$ws.Range('A7').Value = 'test_sql/5'
$ws.Range('B7').Value = 'List all names by customers above the average age | customers: customer_id, name, surname, age'
$ws.Range('C7').Value = '[''SELECT NAME FROM CUSTOMERS WHERE AGE > (SELECT AVG(AGE) FROM CUSTOMERS)'']'
$ws.Range('D7').Value = '[False]'

$ws.Range('A8').Value = 'test_sql/6'
$ws.Range('B8').Value = 'How many orders have products that their prices sum up to 100 or more? | sales: sale_id, product_id, customer_id, quantity, store | orders : order_id, description | order_product : order_id, product_id, price | products : product_id, country, price'
$ws.Range('C8').Value = '[''SELECT COUNT(*) FROM ORDERS AS T1 JOIN ORDER_PRODUCT AS T2 ON T1.ORDER_ID = T2.ORDER_ID WHERE T2.PRICE >= 100'']'
$ws.Range('D8').Value = '[False]'

$ws.Range('A9').Value = 'test_sql/7'
$ws.Range('B9').Value = 'For all sales , what is the most frequent store? | sales: sale_id, product_id, customer_id, quantity, store'
$ws.Range('C9').Value = '[''SELECT STORE, COUNT(STORE) FROM SALES GROUP BY STORE ORDER BY COUNT(STORE) DESC LIMIT 1'']'
$ws.Range('D9').Value = '[False]'

$ws.Range('A10').Value = 'test_sql/8'
$ws.Range('B10').Value = 'Show the sale that has the spanish product with higher quantity | sales: sale_id, product_id, customer_id, quantity, store |  orders : order_id, description | order_product : order_id, product_id, price | products : product_id, country, price'
$ws.Range('C10').Value = '["SELECT T1.SALE_ID FROM SALES AS T1 JOIN ORDER_PRODUCT AS T2 ON T1.PRODUCT_ID = T2.PRODUCT_ID WHERE T2.COUNTRY = ''SPAIN'' AND T1.QUANTITY > 1"]'
$ws.Range('D10').Value = '[False]'

$ws.Range('A11').Value = 'test_sql/9'
$ws.Range('B11').Value = 'Who is the customer that has more sales than the rest? | sales: sale_id, product_id, customer_id, quantity, store | customers: customer_id, name, surname, age'
$ws.Range('C11').Value = '['' SELECT NAME FROM CUSTOMERS AS T1 JOIN SALES AS T2 ON T1.CUSTOMER_ID = T2.CUSTOMER_ID GROUP BY T1.CUSTOMER_ID ORDER BY SUM(T2.SALES_ID) DESC LIMIT 1'']'
$ws.Range('D11').Value = '[False]'

# Header row styling: bold font, thin border all around, centered horizontally, top vertically
$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

